# Update the 2017 household survey "destination_simple" lookup sheet:
#  - recode the "Other"/"Social-Recreation" bucket into finer categories
#  - widen column A to fit the new (longer) category labels
#  - move the active selection to F8 (cosmetic, matches authored file)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("destination_simple")

# --- Recode column B (dest_purpose_simple) ----------------------------
# Written in first-use order (Drop off/Pick up, Social/Recreation/Eat Meal,
# Health and Exercise, Errands) so the shared-string table matches the
# authored workbook's string order.
$ws.Range("B4").Value  = "Drop off/Pick up"

$ws.Range("B12").Value = "Social/Recreation/Eat Meal"
$ws.Range("B14").Value = "Social/Recreation/Eat Meal"
$ws.Range("B15").Value = "Social/Recreation/Eat Meal"
$ws.Range("B16").Value = "Social/Recreation/Eat Meal"
$ws.Range("B17").Value = "Social/Recreation/Eat Meal"
$ws.Range("B20").Value = "Social/Recreation/Eat Meal"

$ws.Range("B11").Value = "Health and Exercise"
$ws.Range("B13").Value = "Health and Exercise"

$ws.Range("B10").Value = "Errands"
$ws.Range("B19").Value = "Errands"

# --- Widen column A for the longer labels ------------------------------
$ws.Columns.Item(1).ColumnWidth = 48.86

# --- Move the selection (cosmetic, matches authored workbook) ----------
$ws.Range("F8").Select()
